# Automatische test-sync: 2025-08-05 18:19:50
# Adds a new log row (#28) to the "Logs" sheet for Testmail #7 and bumps
# the "Inkoop / Bestellingen" tally on the "Dashboard" sheet from 2 to 3.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 28

$logs.Cells.Item($newRow, 1).Value = "Is dit artikel nog op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #7: Is dit artikel nog op voorraad?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-05 18:19:30"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Expand the conditional-formatting ranges (D/G/H/I/J) so they keep
# covering the data through the newly added row.
$colsToExpand = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExpand) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "27")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "28")
    $fcs = $oldRange.FormatConditions
    if ($fcs.Count -gt 0) {
        $fcs.Item(1).ModifyAppliesToRange($newRange)
    }
}

# Bump the Dashboard count for "Inkoop / Bestellingen" (row 4) from 2 to 3.
$dashboard.Range("B4").Value = 3
